# Generate Report for Handoff
#
# The b6121ff3-8af0-4070-aa99-c314affb1e78.md record has been handed off
# again: its Status moves from "Handed back: in sync with en-US" to
# "Ready for handoff" on the Overview sheet (both locale columns) and on
# each per-locale detail sheet, and the handoff timestamps for that run
# are refreshed.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

# Status for the b6121ff3 row (row 3) flips to "Ready for handoff" for
# both the zh-cn (B) and de-de (C) columns.
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# "Latest Handoff Date" is refreshed for both rows (they shared the same
# timestamp text before the edit, and still do after).
$overview.Range("D2").Value = "2016-03-19 04:10:27"
$overview.Range("D3").Value = "2016-03-19 04:10:27"

# ---- zh-cn detail sheet ----------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

# Status for the b6121ff3 row (row 3) flips to "Ready for handoff".
$zhcn.Range("C3").Value = "Ready for handoff"

# "Latest Handoff Datetime" refreshed for both rows.
$zhcn.Range("E2").Value = "2016-03-19 04:10:19"
$zhcn.Range("E3").Value = "2016-03-19 04:10:19"

# ---- de-de detail sheet ----------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

# Status for the b6121ff3 row (row 3) flips to "Ready for handoff".
$dede.Range("C3").Value = "Ready for handoff"

# "Latest Handoff Datetime" refreshed for both rows.
$dede.Range("E2").Value = "2016-03-19 04:10:27"
$dede.Range("E3").Value = "2016-03-19 04:10:27"
